$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 124; this shifts the existing rows 124-191
# down to 125-192, matching the rest of the diff (which is just every row
# from the old 124 onward appearing one row lower).
$ws.Rows("124:124").Insert()

# Populate the newly inserted row 124 with the new weekly data point.
$ws.Cells.Item(124, 1).Value = 3
$ws.Cells.Item(124, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(124, 3).Value = "Coquimbo"
$ws.Cells.Item(124, 4).Value = 44806
$ws.Cells.Item(124, 5).Value = 5
$ws.Cells.Item(124, 6).Value = 100112026
$ws.Cells.Item(124, 7).Value = "Haba"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 105
$ws.Cells.Item(124, 11).Value = 13000
$ws.Cells.Item(124, 12).Value = 14000
$ws.Cells.Item(124, 13).Value = 13524
$ws.Cells.Item(124, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(124, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(124, 16).Value = 541
$ws.Cells.Item(124, 17).Value = 25
$ws.Cells.Item(124, 18).Value = "Hortaliza"
